# Mais testes e classes de equivalencia
#
# Appends four new "List Paragraph" bullet items (numId=1) to the end of
# the document body, right before the final <w:sectPr>, mirroring text
# that was typed in Word at outline levels 2 and 3 (w:ilvl 1 and 2).
#
# We build the four paragraphs as literal WordprocessingML and hand them
# to Range.InsertXML at a collapsed range positioned at the very end of
# the document's main story - this lets us reproduce the exact paragraph
# properties (style/numbering/spacing/justification/run fonts) as well as
# the <w:lastRenderedPageBreak/> marker Word had stamped on the last new
# paragraph, which plain Range.Text/InsertParagraphAfter calls cannot add.

$d = $word.ActiveDocument

$newParagraphsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Após criar uma quantidade de testes, como para leilões, não precisamos testar até 1000 lances para ver se o código continua funcionando.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Isso porque existe algo chamado teste de equivalência, onde temos a confiança de que nosso código funcionará independente do número de lances e se eles foram adicionados em ordem crescente ou decrescente, uma vez que já fizemos esses testes com um número menor de lances.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Portanto, se já fizemos testes com 1, 2, 3 lances, em ordem crescente e decrescente, já podemos confiar que funcionará para inúmeros lances devido ao teste de equivalência. A menos que exista uma condição diferente em algum determinado lance ou que faça sentido ser testado, aí, nesse caso, seria interessante fazer os testes necessários.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Considerando isso é importante sempre levarmos em conta essas classes de equivalência quando formos escrever os testes, para que não percamos tempo escrevendo testes que irão apenas testar mais do mesmo.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Collapsed range sitting right at the end of the document's content
# (before the very final end-of-story mark), so InsertXML appends rather
# than replacing any existing text.
$endOfDoc = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endOfDoc.InsertXML($newParagraphsXml)

Write-Host "Paragraphs in document after edit: $($d.Paragraphs.Count)"
